$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.618
$ws.Range("A3").Value = -21.603
$ws.Range("D5").Value = -7.917999999999999
$ws.Range("E5").Value = 12.94
$ws.Range("E9").Value = 13.321
$ws.Range("E11").Value = 13.132
$ws.Range("A14").Value = -20.814
$ws.Range("A16").Value = -20.591
$ws.Range("D16").Value = -8.161
$ws.Range("E17").Value = 13.784
$ws.Range("A21").Value = -21.04
$ws.Range("E21").Value = 13.44
$ws.Range("A23").Value = -21.584
$ws.Range("A25").Value = -22.27
